$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") values for rows 2-10 move from 45183 (2023-09-14)
# to 45184 (2023-09-15) while keeping their existing date formatting/style.
for ($row = 2; $row -le 10; $row++) {
    $cell = $ws.Range("C$row")
    if ($cell.Value2 -eq 45183) {
        $cell.Value2 = 45184
    }
}
